# Fruta / hortaliza, semanal
# Insert a new data row at row 49 (pushing the existing rows 49:145 down to
# 50:146) and populate it with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 49; this shifts rows 49:145 down to
# 50:146 (and correspondingly grows the sheet dimension to A1:R146).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44883
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 100114002
$ws.Range("G49").Value = "Camote"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 30
$ws.Range("K49").Value = 24000
$ws.Range("L49").Value = 24000
$ws.Range("M49").Value = 24000
$ws.Range("N49").Value = "$/malla 20 kilos"
$ws.Range("O49").Value = "Perú"
$ws.Range("P49").Value = 1200
$ws.Range("Q49").Value = 20
$ws.Range("R49").Value = "Hortaliza"
